$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add saving original clip data to columns C and D, rows 1-3
$ws.Range("C1").Value = 5244.444444444444
$ws.Range("D1").Value = "CK1_1_0002.jpg"

$ws.Range("C2").Value = 1662.698412698413
$ws.Range("D2").Value = "CK1_1_0002.jpg"

$ws.Range("C3").Value = 7256.349206349206
$ws.Range("D3").Value = "CK1_1_0002.jpg"
